$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# SignInTest (sheet1): swap in the new credentials, drop the mailto: hyperlinks
# ---------------------------------------------------------------------------
$ws1.Hyperlinks.Delete()
$ws1.Range("A2:B3").Style = "Normal"

$ws1.Range("A2").Value = "ramesh.sanjay@gmail.com"
$ws1.Range("B2").Value = "Selenium@1234"
$ws1.Range("C2").Value = "Y"

$ws1.Range("A3").Value = "uip311977@gmail.com"
$ws1.Range("B3").Value = "Selenium@1234"
$ws1.Range("C3").Value = "Y"

# Stray formatted cell that shows up in the authored workbook.
$ws1.Range("G21").Style = "Hyperlink"
$ws1.Range("G21").ClearContents()

$ws1.Columns.Item(1).ColumnWidth = 25.307291666666668
$ws1.Columns.Item(2).ColumnWidth = 13.592447916666666

$ws1.Range("C10").Select()

# ---------------------------------------------------------------------------
# FlightSearchTest (sheet2): new route/dates, drop noOfAdults/noOfChildern
# ---------------------------------------------------------------------------
$ws2.Columns.Item(5).Delete()
$ws2.Columns.Item(5).Delete()

$ws2.Range("C2").Copy()
$ws2.Range("D2").PasteSpecial(-4122)

$ws2.Range("A2").Value = "Bengaluru, India (BLR-Kempegowda Intl.)"
$ws2.Range("B2").Value = "Goa, India (GOI-Dabolim)"
$ws2.Range("C2").Value = "20/2/2020"
$ws2.Range("D2").Value = "26/2/2020"
$ws2.Range("E2").Value = "Y"

$ws2.Columns.Item(1).ColumnWidth = 32.592447916666664
$ws2.Columns.Item(2).ColumnWidth = 49.166666666666664

$ws2.Range("A7").Select()
$ws2.Activate()
